$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new state column (before column B) and a new state row (before row 3) ---
# This mirrors inserting an extra "0" state into the s_n/x_n table: existing columns B:I
# shift right to C:J, existing rows 3:8 shift down to 4:9.
$ws.Columns("B:B").Insert()
$ws.Rows("3:3").Insert()

# The freshly inserted column/row picked up formatting from their left/above neighbour
# (column A / row 2 respectively). Re-pull the correct formatting from the matching
# original column/row (now shifted one to the right / down) so the new column & row look
# like the rest of the table.
$ws.Range("C2:C9").Copy()
$ws.Range("B2:B9").PasteSpecial(-4122)

$ws.Range("A4:J4").Copy()
$ws.Range("A3:J3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Re-establish the merged header cell across the widened header row ---
# The header text ("f_n(s_n,x_n)") lives in the old merge anchor, which is now one column
# to the right (C1) after the column insert; move it back to the new anchor (B1) before
# merging, since merging keeps only the top-left cell's content.
$txt = $ws.Cells.Item(1, 3).Value()
$ws.Cells.Item(1, 2).Value = $txt
$ws.Cells.Item(1, 3).Value = ""

# Stash the header's fill/border/alignment formatting (from D1, still a plain cell at this
# point) in a scratch cell so it survives the merge operation below.
$ws.Range("D1:D1").Copy()
$ws.Range("Z1:Z1").PasteSpecial(-4122)

$ws.Range("B1:H1").Merge()

# Re-apply the stashed header formatting across the whole merged header cell, then clean up
# the scratch cell.
$ws.Range("Z1:Z1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$ws.Range("Z1:Z1").Clear()

$excel.CutCopyMode = 0

# --- Fill in the values for the newly inserted "0" state column & row ---
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(3, 1).Value = 0

# --- Restore the saved selection/cursor position ---
$ws.Range("E15").Select()
